# "Conversions work, formatting tidied up a bit"
# - Remove the "Variations" column (C) contents: the header (C1) and the
#   three stray variation notes (C30, C37, C43). Column D ("Unmarked")
#   stays put, it is not shifted left.
# - Move the active selection from C25 to D15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C1:C44").ClearContents()

$ws.Range("D15").Select() | Out-Null
